$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.57"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.862.84"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'304.84"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.5050"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("D8").Value = "'0.3642"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").Value = "'0.07163"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'0.8927"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'20.70"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.869.75"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'0.07479"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'92.56"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "'5.228"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.000008492"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'14.20"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "26.940.03"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "'5.028"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "2.093.25"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'10.38"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").Value = "'6.398"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'146.95"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").Value = "'1.791"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").Value = "'17.87"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "'2.080"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "'113.15"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'4.704"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").Value = "'4.674"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'0.09250"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").Value = "'0.05089"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'0.7515"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").Value = "'2.996"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("D36").Value = "'1.151"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "'3.268"
$ws.Range("E37").Value = "  +6.58%  "
$ws.Range("D38").Value = "'2.534"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "'0.01997"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("D41").Value = "'1.071"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'118.54"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'6.534"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'8.525"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'0.4688"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'10.08"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'62.90"
$ws.Range("E51").Value = "  -2.27%  "
